$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark that currently sits
#    right after "It's" (before " simple yet clear...").
# ---------------------------------------------------------------
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

# ---------------------------------------------------------------
# 2) Rebuild the very first paragraph ("Name: Quwaine") so the
#    stray spell-check markers (<w:proofErr .../>) around
#    "Quwaine" are gone and " " + "Quwaine" become one run.
#    (The only reliable way to drop the orphan proofErr markers
#    in this engine is to remove the whole paragraph and retype
#    its text.)
# ---------------------------------------------------------------
$firstPara = $d.Range(0, 14)
$firstPara.Delete()
$p1 = $d.Paragraphs(1)
$p1.Range.InsertBefore("Name: Quwaine`r")

# ---------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark at the very start of the
#    document (zero-length, before "Name:"). Adding a bookmark on
#    a plain (0,0) range gets mis-handled by this engine, so we
#    insert a throw-away placeholder character, bookmark that
#    single character, then clear the bookmarked text back out --
#    this leaves a correctly collapsed bookmark sitting at
#    position 0.
# ---------------------------------------------------------------
$placeholder = $d.Range(0, 0)
$placeholder.InsertBefore("X")
$placeholderRange = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange)
$newBm = $d.Bookmarks("_GoBack")
$newBm.Range.Text = ""

# ---------------------------------------------------------------
# 4) Split "Name: Quwaine" back into "Name:" + " Quwaine" runs
#    (matching the original two-run layout) by nudging formatting
#    on the " Quwaine" span and immediately reverting it -- this
#    forces a run boundary without changing any visible formatting.
# ---------------------------------------------------------------
$quwaineRange = $d.Range(5, 13)
$quwaineRange.Bold = 1
$quwaineRange.Bold = 0
